# Auto-update GitHub repos Excel export
# Inserts three new repos (from j-chaganti) at the top of the list and
# shifts the existing Josmietha repos down, dropping the last one
# (repo-scanner) off the bottom of the moved block since it becomes one
# of the newly-inserted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Git_Merge_Conflict_Resolution_Demo", "https://github.com/j-chaganti/Git_Merge_Conflict_Resolution_Demo"),
    @("my-repo", "https://github.com/j-chaganti/my-repo"),
    @("repo-scanner", "https://github.com/j-chaganti/repo-scanner"),
    @("GenAI-Hackathon", "https://github.com/Josmietha/GenAI-Hackathon"),
    @("merge-conflict-exercise", "https://github.com/Josmietha/merge-conflict-exercise"),
    @("merge-conflict-practice", "https://github.com/Josmietha/merge-conflict-practice"),
    @("My-HTML-Portfolio", "https://github.com/Josmietha/My-HTML-Portfolio"),
    @("OCT_Task1", "https://github.com/Josmietha/OCT_Task1"),
    @("OIBSIP", "https://github.com/Josmietha/OIBSIP"),
    @("Practice", "https://github.com/Josmietha/Practice"),
    @("practicee", "https://github.com/Josmietha/practicee")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
